$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "sexo" column (E) dimension was re-curated from a dimension to a measure.
# Update the concept-type cell (E2) and its associated metadata cells (E3, E4).
$ws.Range("E2").Value = "iaest-measure:sexo"
$ws.Range("E3").Value = "medida"
$ws.Range("E4").Value = "xsd:int"

# Row 5 (the skos:Concept / mapping-sexo.xlsx leftovers from the old
# dimension-based mapping) is no longer needed and is removed entirely.
$ws.Range("E5").EntireRow.Delete()
